$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D12").Value = "“머신러닝, 딥러닝 어떻게 공부할 것인가?” 무료 강의 안내"
$ws.Range("E12").Value = "https://tensorflow.blog/2023/10/08/%eb%a8%b8%ec%8b%a0%eb%9f%ac%eb%8b%9d-%eb%94%a5%eb%9f%ac%eb%8b%9d-%ec%96%b4%eb%96%bb%ea%b2%8c-%ea%b3%b5%eb%b6%80%ed%95%a0-%ea%b2%83%ec%9d%b8%ea%b0%80-%eb%ac%b4%eb%a3%8c-%ea%b0%95%ec%9d%98/"

$ws.Range("D44").Value = "[Book Review] 칩워 (Chipwar)"
$ws.Range("E44").Value = "https://engineering-ladder.tistory.com/134"

$ws.Range("D51").Value = "[pandas] 데이터프레임 짝수행, 홀수행 선택하기"
$ws.Range("E51").Value = "https://bskyvision.com/entry/pandas-%EB%8D%B0%EC%9D%B4%ED%84%B0%ED%94%84%EB%A0%88%EC%9E%84-%EC%A7%9D%EC%88%98%ED%96%89-%ED%99%80%EC%88%98%ED%96%89-%EC%84%A0%ED%83%9D%ED%95%98%EA%B8%B0"
